$p = $ppt.ActivePresentation

# --- 1. Update the cached date field text on the slide master and every slide layout ---
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "06/03/2020") {
            $sh.TextFrame.TextRange.Text = "13/03/2020"
        }
    }
}

for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $cl = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "06/03/2020") {
                $sh.TextFrame.TextRange.Text = "13/03/2020"
            }
        }
    }
}

# --- 2. Add the new bullet points to slide 5's content placeholder ---
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(4)
$tf = $shape.TextFrame
$tf.TextRange.Text = "Talk about games : - Gamification`rTal about ML:`rSubtopic of Artificial intelligence. "

$tf.TextRange.Paragraphs(1).Font.Size = 24
$tf.TextRange.Paragraphs(2).Font.Size = 24
$tf.TextRange.Paragraphs(3).Font.Size = 20
$tf.TextRange.Paragraphs(3).IndentLevel = 2
